$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue list")

# --- Insert 6 new rows right after row 15 (the old last blank 'No=11' row), ---
# --- pushing the trailing blank spacer rows further down.                  ---
$ws.Range("16:21").Insert()

# --- Row 14 (No=10): fill in the issue that used to be an empty placeholder ---
$ws.Range("C14").Value = 43458
$ws.Range("D14").Value = "Closed"
$ws.Range("F14").Value = "Power Off가 안됨."
$ws.Range("G14").Value = "KEY_PWR_Pin GPIO 설정 변경`nNOPULL -> PULLDOWN (Gpio.c)"

# --- Row 15 (No=11): fill in the issue that used to be an empty placeholder ---
$ws.Range("C15").Value = 43459
$ws.Range("D15").Value = "Closed"
$ws.Range("F15").Value = "Plasma On Key 동작 안됨.`n Default Low 유지"
$ws.Range("G15").Value = "충전중 Plasma 동작 안하도록 수정하며 생긴 문제임`nUSB_DET_Pin GPIO 설정 변경`nNOPULL -> PULLDOWN (Gpio.c)"

# --- Row 16 (No=12): brand-new issue row ---
$ws.Range("B16").Value = 12
$ws.Range("C16").Value = 43466
$ws.Range("D16").Value = "OPEN"
$ws.Range("F16").Value = "완충 인식 안됨"

# --- Row 17 (No=13): brand-new issue row ---
$ws.Range("B17").Value = 13
$ws.Range("C17").Value = 43481
$ws.Range("D17").Value = "Closed"
$ws.Range("E17").Value = 43482
$ws.Range("F17").Value = "Plasma On후 5초후 S/V On 지원"
$ws.Range("G17").Value = "Plasma_state.c - Line 24 수정`n                GPIO_ENABLE(GAS_EN);`n                vTaskDelay( 5000); // 2019.01.17 Arvid - Plasma On delay 5sec"

# --- Rows 18-21: new blank placeholder issue rows 14-17 ---
$ws.Range("B18").Value = 14
$ws.Range("B19").Value = 15
$ws.Range("B20").Value = 16
$ws.Range("B21").Value = 17

Write-Output "done"
